# "Change One Rev 1"
# Updates the Test Results sheet: several rows of test-case results are
# corrected/re-ordered, two rows are added (CL-29, CT-44), three rows are
# removed (CT-28, CT-71, CT-96) and the last row (CT-96 -> CT-100) is
# refreshed, while some statuses flip between Passed/Failed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact number-format codes already used in the workbook so we reuse the
# existing style entries instead of creating new ones.
$fmtStart = "[h]:mm:ss;@"
$fmtEnd   = "m/d/yy\ h:mm;@"

function Set-ResultRow {
    param($Row, $TestId, $Start, $End, $Status, $Output, $StatusStyle, $EndFormat)

    $ws.Range("B$Row").Value = $TestId

    $ws.Range("C$Row").NumberFormat = $fmtStart
    $ws.Range("C$Row").Value = $Start

    $effectiveEndFormat = $EndFormat
    if (-not $effectiveEndFormat) {
        $effectiveEndFormat = $fmtEnd
    }
    $ws.Range("D$Row").NumberFormat = $effectiveEndFormat
    $ws.Range("D$Row").Value = $End

    $effectiveStyle = $StatusStyle
    if (-not $effectiveStyle) {
        if ($Status -eq "Passed") {
            $effectiveStyle = "Good"
        } else {
            $effectiveStyle = "Bad"
        }
    }
    $ws.Range("E$Row").Style = $effectiveStyle
    $ws.Range("E$Row").Value = $Status

    if ($Output) {
        $ws.Range("F$Row").Value = $Output
    } else {
        $ws.Range("F$Row").ClearContents()
    }
}

# Prime the brand-new shared strings in the same order the original author's
# Excel session produced them in (newest-edited row first), so the saved
# sharedStrings.xml table lines up exactly with the source workbook.
$ws.Range("B23").Value = "CT-100"
$ws.Range("B14").Value = "CT-44"
$ws.Range("B12").Value = "CL-29"

Set-ResultRow 4  "CT-10" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:33 PM" "Passed" $null
Set-ResultRow 5  "CT-14" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:29 PM" "Passed" $null
Set-ResultRow 6  "CT-22" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:33 PM" "Passed" $null
Set-ResultRow 7  "CT-23" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:33 PM" "Failed" $null
Set-ResultRow 8  "CT-24" "9/12/2017  13:50:30 PM" "9/12/2017  13:50:57 PM" "Passed" $null
Set-ResultRow 9  "CT-25" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:29 PM" "Passed" $null
Set-ResultRow 10 "CT-26" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:29 PM" "Passed" $null
Set-ResultRow 11 "CT-27" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:29 PM" "Failed" "Exception: object not found"
Set-ResultRow 12 "CL-29" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:29 PM" "Passed" $null
Set-ResultRow 13 "CT-45" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:27 PM" "Passed" $null

Set-ResultRow 14 "CT-44" 44540.448611111111 44540.490277777797 "Failed" $null "Good" $fmtStart

Set-ResultRow 15 "CT-46" "9/12/2017  13:50:36 PM" "9/12/2017  13:50:36 PM" "Passed" $null
Set-ResultRow 16 "CT-58" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:27 PM" "Failed" $null
Set-ResultRow 17 "CT-62" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:27 PM" "Passed" $null
Set-ResultRow 18 "CT-65" "9/12/2017  13:50:28 PM" "9/12/2017  13:50:28 PM" "Passed" $null
Set-ResultRow 19 "CT-66" "9/12/2017  13:50:29 PM" "9/12/2017  13:50:29 PM" "Passed" $null
Set-ResultRow 20 "CT-73" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:27 PM" "Failed" "Connection could not be established"
Set-ResultRow 21 "CT-73" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:27 PM" "Passed" $null
Set-ResultRow 22 "CT-74" "9/12/2017  13:50:27 PM" "9/12/2017  13:50:27 PM" "Failed" $null
Set-ResultRow 23 "CT-100" "9/12/2017  13:50:27 PM" "9/12/2017  14:02:14 PM" "Passed" $null

$ws.Range("B12").Select()
